$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$val_B2 = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
WHERE diag.stage_of_disease IN ['Vb']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@
$ws.Range("B2").Value = $val_B2

$val_B3 = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
 WHERE diag.stage_of_disease IN ['Vb']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc
limit 100
'@
$ws.Range("B3").Value = $val_B3

$val_C3 = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.stage_of_disease IN ['Vb']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`+D3
'@
$ws.Range("C3").Value = $val_C3

$val_B4 = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE  diag.stage_of_disease in ['Vb']
OPTIONAL MATCH (f)-[*]->(samp:sample)
WITH
        DISTINCT f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
        order by f.file_name asc
        limit 100
'@
$ws.Range("B4").Value = $val_B4

$val_B5 = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (c)<--(demo:demographic)
WHERE  diag.stage_of_disease in ['Vb']
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@
$ws.Range("B5").Value = $val_B5

$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 261
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 377

$styledRange = $ws.Range("B2:C5")
$styledRange.Interior.ThemeColor = 10
$styledRange.Interior.TintAndShade = 0.59999389629810485
$styledRange.Borders.LineStyle = 1
$styledRange.Borders.Weight = 2

$ws.Range("D5:G5").Select()
$excel.ActiveWindow.Zoom = 82

"Edit complete"